# Update gh-pages to output generated at 456a3b4
# Applies F/G column updates (想去人数 / 最低票价) to sheets "展览" (1) and "全部类型" (4)

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F7").Value = 14988
$ws1.Range("F10").Value = 679
$ws1.Range("F11").Value = 15234
$ws1.Range("F13").Value = 8758
$ws1.Range("F16").Value = 66
$ws1.Range("F20").Value = 13
$ws1.Range("F22").Value = 514
$ws1.Range("F32").Value = 417
$ws1.Range("G32").Value = "不可售"
$ws1.Range("F33").Value = 32
$ws1.Range("F35").Value = 230
$ws1.Range("F38").Value = 111
$ws1.Range("F39").Value = 5376

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F7").Value = 14988
$ws4.Range("F10").Value = 679
$ws4.Range("F11").Value = 15234
$ws4.Range("F13").Value = 8758
$ws4.Range("F17").Value = 66
$ws4.Range("F21").Value = 13
$ws4.Range("F23").Value = 514
$ws4.Range("F35").Value = 417
$ws4.Range("G35").Value = "不可售"
$ws4.Range("F36").Value = 32
$ws4.Range("F38").Value = 230
$ws4.Range("F41").Value = 111
$ws4.Range("F42").Value = 5376
